$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename hospital entries to new Frisius MC locations (merger/renaming of
# Medisch Centrum Leeuwarden and Tjongerschans into Frisius MC)
$ws.Range("A33").Value = "Frisius MC locatie Leeuwarden"
$ws.Range("A63").Value = "Frisius MC locatie Heerenveen"
